# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (D5) and "Correspond Handback DateTime" (G5)
# timestamps on the zh-cn and de-de sheets to reflect the new handback report run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-18 08:23:01"
$wsZhCn.Range("G5").Value = "2016-02-18 08:23:44"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-18 08:23:15"
$wsDeDe.Range("G5").Value = "2016-02-18 08:24:07"
